# Weekly update: insert a new price record row at row 75 (Choclo, "Dulce o
# Americano", Provincia de Limarí) pushing the existing rows 75-102 down to
# 76-103, and extending the sheet's used range accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 75; this shifts rows 75..102 down to
# 76..103 and copies formatting (e.g. the date style on column D) from the
# row that used to occupy position 75.
$ws.Rows.Item(75).Insert()

# Populate the newly inserted row 75 with the new record.
$ws.Range("A75").Value = 7
$ws.Range("B75").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C75").Value = "Ñuble"
$ws.Range("D75").Value = 44523
$ws.Range("E75").Value = 16
$ws.Range("F75").Value = 100112024
$ws.Range("G75").Value = "Choclo"
$ws.Range("H75").Value = "Dulce o Americano"
$ws.Range("I75").Value = "Primera"
$ws.Range("J75").Value = 80
$ws.Range("K75").Value = 16000
$ws.Range("L75").Value = 17000
$ws.Range("M75").Value = 16500
$ws.Range("N75").Value = "`$/malla 60 unidades"
$ws.Range("O75").Value = "Provincia de Limarí"
$ws.Range("P75").Value = 275
$ws.Range("Q75").Value = 60
$ws.Range("R75").Value = "Hortaliza"
